$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr2 = $sh.TextFrame2.TextRange
$para = $tr2.Paragraphs(2)
$para.Text = "Requisitos: El grupo debe entregar previamente todas las actividades y ejercicios planteados."
